# Auto-generated edit script: updates market-price-derived columns (H-N)
# on the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to match the refreshed
# Universalis price snapshot pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 914.5
$ws.Range("I2").Value = 803
$ws.Range("K2").Value = 803
$ws.Range("M2").Value = -690
$ws.Range("H15").Value = 1018.02325
$ws.Range("I15").Value = 1018.02325
$ws.Range("K15").Value = 3054.06975
$ws.Range("M15").Value = -2885.06975
$ws.Range("H17").Value = 1732.6364
$ws.Range("J17").Value = 1799.6
$ws.Range("L17").Value = 5398.799999999999
$ws.Range("N17").Value = -5734.799999999999
$ws.Range("H28").Value = 288.4
$ws.Range("I28").Value = 288.4
$ws.Range("K28").Value = 288.4
$ws.Range("M28").Value = 196.6
$ws.Range("H55").Value = 681.2
$ws.Range("I55").Value = 613.3333
$ws.Range("J55").Value = 783
$ws.Range("K55").Value = 613.3333
$ws.Range("L55").Value = 783
$ws.Range("M55").Value = -399.3333
$ws.Range("N55").Value = -1211
$ws.Range("H70").Value = 97846.375
$ws.Range("I70").Value = 2999.5
$ws.Range("J70").Value = 129462
$ws.Range("K70").Value = 8998.5
$ws.Range("L70").Value = 388386
$ws.Range("M70").Value = -8728.5
$ws.Range("N70").Value = -388926
$ws.Range("H73").Value = 97846.375
$ws.Range("I73").Value = 2999.5
$ws.Range("J73").Value = 129462
$ws.Range("K73").Value = 8998.5
$ws.Range("L73").Value = 388386
$ws.Range("M73").Value = -8062.5
$ws.Range("N73").Value = -390258
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H112").Value = 1864.1765
$ws.Range("J112").Value = 2118.6155
$ws.Range("L112").Value = 6355.8465
$ws.Range("N112").Value = -8571.8465
$ws.Range("H132").Value = 1478.15
$ws.Range("I132").Value = 1478.15
$ws.Range("K132").Value = 4434.450000000001
$ws.Range("M132").Value = -1904.450000000001
$ws.Range("H135").Value = 1098.8235
$ws.Range("I135").Value = 727.1429000000001
$ws.Range("K135").Value = 6544.2861
$ws.Range("M135").Value = -4009.2861
$ws.Range("H141").Value = 3087.7917
$ws.Range("J141").Value = 11375.5
$ws.Range("L141").Value = 34126.5
$ws.Range("N141").Value = -44486.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7654.1904
$ws.Range("I32").Value = 5482.4
$ws.Range("K32").Value = 5482.4
$ws.Range("M32").Value = -5195.4
$ws.Range("H74").Value = 770.4583
$ws.Range("I74").Value = 770.4583
$ws.Range("K74").Value = 770.4583
$ws.Range("M74").Value = 103.5417
$ws.Range("H77").Value = 770.4583
$ws.Range("I77").Value = 770.4583
$ws.Range("K77").Value = 3852.2915
$ws.Range("M77").Value = 515.7084999999997
$ws.Range("H110").Value = 2694.625
$ws.Range("I110").Value = 2208.2666
$ws.Range("K110").Value = 2208.2666
$ws.Range("M110").Value = -163.2665999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2225.4736
$ws.Range("I105").Value = 1958
$ws.Range("J105").Value = 4499
$ws.Range("K105").Value = 1958
$ws.Range("L105").Value = 4499
$ws.Range("M105").Value = -211
$ws.Range("N105").Value = -7993
$ws.Range("H134").Value = 1409.7715
$ws.Range("I134").Value = 1047.2
$ws.Range("J134").Value = 1893.2
$ws.Range("K134").Value = 3141.6
$ws.Range("L134").Value = 5679.6
$ws.Range("M134").Value = -606.6000000000004
$ws.Range("N134").Value = -10749.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 447.66666
$ws.Range("I22").Value = 447
$ws.Range("K22").Value = 447
$ws.Range("M22").Value = -97
$ws.Range("H58").Value = 2544.8462
$ws.Range("I58").Value = 1316.5555
$ws.Range("K58").Value = 1316.5555
$ws.Range("M58").Value = -1113.5555
$ws.Range("H132").Value = 2577.1333
$ws.Range("I132").Value = 2381.3845
$ws.Range("K132").Value = 7144.1535
$ws.Range("M132").Value = -4614.1535
$ws.Range("H134").Value = 3020.7273
$ws.Range("I134").Value = 2998.2144
$ws.Range("K134").Value = 8994.643199999999
$ws.Range("M134").Value = -6459.643199999999
$ws.Range("H136").Value = 2544.8462
$ws.Range("I136").Value = 1316.5555
$ws.Range("K136").Value = 3949.6665
$ws.Range("M136").Value = -1399.6665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1616792.2
$ws.Range("J4").Value = 2710.2
$ws.Range("L4").Value = 8130.599999999999
$ws.Range("N4").Value = -8354.599999999999
$ws.Range("H5").Value = 954.63635
$ws.Range("I5").Value = 541.6
$ws.Range("J5").Value = 1298.8334
$ws.Range("K5").Value = 1624.8
$ws.Range("L5").Value = 3896.5002
$ws.Range("M5").Value = -1512.8
$ws.Range("N5").Value = -4120.5002
$ws.Range("H38").Value = 7279.357
$ws.Range("J38").Value = 81.40000000000001
$ws.Range("L38").Value = 244.2
$ws.Range("N38").Value = -938.2
$ws.Range("H132").Value = 3075.3914
$ws.Range("I132").Value = 3504.2856
$ws.Range("J132").Value = 2408.2222
$ws.Range("K132").Value = 31538.5704
$ws.Range("L132").Value = 21673.9998
$ws.Range("M132").Value = -29008.5704
$ws.Range("N132").Value = -26733.9998
$ws.Range("H135").Value = 954.63635
$ws.Range("I135").Value = 541.6
$ws.Range("J135").Value = 1298.8334
$ws.Range("K135").Value = 4874.400000000001
$ws.Range("L135").Value = 11689.5006
$ws.Range("M135").Value = -2339.400000000001
$ws.Range("N135").Value = -16759.5006
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("M138").ClearContents()
$ws.Range("N138").ClearContents()
$ws.Range("H140").Value = 5000
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6624.625
$ws.Range("I70").Value = 5249.5
$ws.Range("J70").Value = 7999.75
$ws.Range("K70").Value = 5249.5
$ws.Range("L70").Value = 7999.75
$ws.Range("M70").Value = -4979.5
$ws.Range("N70").Value = -8539.75
$ws.Range("H73").Value = 6624.625
$ws.Range("I73").Value = 5249.5
$ws.Range("J73").Value = 7999.75
$ws.Range("K73").Value = 5249.5
$ws.Range("L73").Value = 7999.75
$ws.Range("M73").Value = -4313.5
$ws.Range("N73").Value = -9871.75
$ws.Range("H80").Value = 5160.4287
$ws.Range("I80").Value = 2999.6667
$ws.Range("J80").Value = 6781
$ws.Range("K80").Value = 2999.6667
$ws.Range("L80").Value = 6781
$ws.Range("M80").Value = -2001.6667
$ws.Range("N80").Value = -8777
$ws.Range("H83").Value = 5160.4287
$ws.Range("I83").Value = 2999.6667
$ws.Range("J83").Value = 6781
$ws.Range("K83").Value = 14998.3335
$ws.Range("L83").Value = 33905
$ws.Range("M83").Value = -10006.3335
$ws.Range("N83").Value = -43889
$ws.Range("H113").Value = 3665.1333
$ws.Range("I113").Value = 2498.375
$ws.Range("K113").Value = 2498.375
$ws.Range("M113").Value = -328.375
$ws.Range("H129").Value = 49999.5
$ws.Range("J129").Value = 49999.5
$ws.Range("L129").Value = 49999.5
$ws.Range("N129").Value = -59999.5
$ws.Range("H137").Value = 16709
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2098.875
$ws.Range("I7").Value = 1514.3846
$ws.Range("K7").Value = 1514.3846
$ws.Range("M7").Value = -1402.3846
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H40").Value = 3789
$ws.Range("I40").Value = 3789
$ws.Range("K40").Value = 3789
$ws.Range("M40").Value = -3653
$ws.Range("H43").Value = 1009599.4
$ws.Range("J43").Value = 1259999.5
$ws.Range("L43").Value = 1259999.5
$ws.Range("N43").Value = -1260385.5
$ws.Range("H61").Value = 2798.8
$ws.Range("I61").Value = 2554.2222
$ws.Range("K61").Value = 2554.2222
$ws.Range("M61").Value = -2352.2222
$ws.Range("H68").Value = 1925.3334
$ws.Range("J68").Value = 2567.3333
$ws.Range("L68").Value = 2567.3333
$ws.Range("N68").Value = -4065.3333
$ws.Range("H71").Value = 1925.3334
$ws.Range("J71").Value = 2567.3333
$ws.Range("L71").Value = 12836.6665
$ws.Range("N71").Value = -20324.6665
$ws.Range("H93").Value = 1273.4615
$ws.Range("I93").Value = 1121.75
$ws.Range("K93").Value = 1121.75
$ws.Range("M93").Value = 126.25
$ws.Range("H106").Value = 17330.334
$ws.Range("J106").Value = 17330.334
$ws.Range("L106").Value = 17330.334
$ws.Range("N106").Value = -19854.334
$ws.Range("H113").Value = 2798.8
$ws.Range("I113").Value = 2554.2222
$ws.Range("K113").Value = 2554.2222
$ws.Range("M113").Value = -384.2222000000002
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H126").Value = 2098.875
$ws.Range("I126").Value = 1514.3846
$ws.Range("K126").Value = 4543.1538
$ws.Range("M126").Value = -2073.1538
$ws.Range("H136").Value = 8844
$ws.Range("J136").Value = 8844
$ws.Range("L136").Value = 26532
$ws.Range("N136").Value = -31632

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 44678.8
$ws.Range("J64").Value = 44678.8
$ws.Range("L64").Value = 44678.8
$ws.Range("N64").Value = -45174.8
$ws.Range("H67").Value = 44678.8
$ws.Range("J67").Value = 44678.8
$ws.Range("L67").Value = 44678.8
$ws.Range("N67").Value = -46394.8
$ws.Range("H75").Value = 60000
$ws.Range("J75").Value = 60000
$ws.Range("L75").Value = 60000
$ws.Range("N75").Value = -61872
$ws.Range("H78").Value = 60000
$ws.Range("J78").Value = 60000
$ws.Range("L78").Value = 180000
$ws.Range("N78").Value = -189360
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H113").Value = 867.4545000000001
$ws.Range("I113").Value = 849.6667
$ws.Range("J113").Value = 888.8
$ws.Range("K113").Value = 2549.0001
$ws.Range("L113").Value = 2666.4
$ws.Range("M113").Value = -379.0001000000002
$ws.Range("N113").Value = -7006.4
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H126").Value = 2198.75
$ws.Range("I126").Value = 2284.6428
$ws.Range("J126").Value = 1597.5
$ws.Range("K126").Value = 6853.928400000001
$ws.Range("L126").Value = 4792.5
$ws.Range("M126").Value = -4383.928400000001
$ws.Range("N126").Value = -9732.5
$ws.Range("H136").Value = 1330.1364
$ws.Range("I136").Value = 914.6111
$ws.Range("K136").Value = 2743.8333
$ws.Range("M136").Value = -193.8332999999998

Write-Output "Applied Seraph_Profits price refresh."